# aggiornamento a 9/09 compreso
# Append new daily rows (367-374, dates 2021-09-02 .. 2021-09-09) to Sheet1,
# matching the style of the existing data (column A formatted as a date,
# style copied from the last existing row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (366) down into the
# new rows so the new date cells (column A) pick up the same style (s="2")
# instead of Excel's default formatting.
$ws.Range("A366").Copy()
$ws.Range("A367:A374").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data, one row per day through 2021-09-09 inclusive.
$data = @(
    @(367, 44441, 0, 5, 156.2988433885589),
    @(368, 44442, 0, 5, 156.2988433885589),
    @(369, 44443, 0, 5, 156.2988433885589),
    @(370, 44444, 0, 0, 0),
    @(371, 44445, 0, 0, 0),
    @(372, 44446, 0, 0, 0),
    @(373, 44447, 0, 0, 0),
    @(374, 44448, 1, 1, 31.25976867771178)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 1).Value = $entry[1]
    $ws.Cells.Item($row, 2).Value = $entry[2]
    $ws.Cells.Item($row, 3).Value = $entry[3]
    $ws.Cells.Item($row, 4).Value = $entry[4]
}
